# Weekly refresh of Fruta / Hortaliza data: the price/volume records for
# rows 2-21 get reshuffled among themselves (dates, quality grade, volume,
# min/max/avg price and $/Kg all travel together as one record). Rows 18,
# 19 and 21 keep their original data; the rest are redistributed per the
# mapping below (new row -> source/original row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (source row as it existed BEFORE this edit)
$mapping = @{
    2  = 4
    3  = 17
    4  = 14
    5  = 13
    6  = 8
    7  = 20
    8  = 6
    9  = 7
    10 = 15
    11 = 5
    12 = 9
    13 = 2
    14 = 3
    15 = 10
    16 = 11
    17 = 12
    18 = 18
    19 = 19
    20 = 16
    21 = 21
}

# Snapshot the current (pre-edit) values for every row first, so that
# writes to earlier rows don't clobber data still needed for later rows.
$snapshot = @{}
for ($r = 2; $r -le 21; $r++) {
    $row = @{
        D = $ws.Cells.Item($r, 4).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
    $snapshot[$r] = $row
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value = $data.D
    $ws.Cells.Item($destRow, 9).Value = $data.I
    $ws.Cells.Item($destRow, 10).Value = $data.J
    $ws.Cells.Item($destRow, 11).Value = $data.K
    $ws.Cells.Item($destRow, 12).Value = $data.L
    $ws.Cells.Item($destRow, 13).Value = $data.M
    $ws.Cells.Item($destRow, 16).Value = $data.P
}
